# The document ends with a run of 12 empty paragraphs (right after the
# "git config --global user.email ..." paragraph, just before the
# section properties). The edit trims 8 of those empty paragraphs,
# leaving 4 trailing empty paragraphs.

$d = $word.ActiveDocument

$keep = 4
$toDelete = 8

for ($i = 0; $i -lt $toDelete; $i++) {
    $total = $d.Paragraphs.Count
    $idx = $total - $keep
    $p = $d.Paragraphs.Item($idx)
    $p.Range.Delete()
}
